# Rename the two logo pictures that live in the document's first-page /
# default headers and footers.
#
#   - Pearson Edexcel logo (footers)   : image2.png -> image1.png
#   - BTEC logo            (headers)   : image1.jpg -> image2.jpg
#
# The pictures themselves are untouched; only the shape's display Name
# (OOXML <wp:docPr name="...">) changes, matching the source diff.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: Pearson Edexcel logo, image2.png -> image1.png ---------------
for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists -and $ftr.Range.ShapeRange.Count -gt 0) {
        $shp = $ftr.Range.ShapeRange.Item(1)
        if ($shp.Name -eq "image2.png") {
            $shp.Name = "image1.png"
        }
    }
}

# --- Headers: BTEC logo, image1.jpg -> image2.jpg ---------------------------
for ($h = 1; $h -le 2; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists -and $hdr.Range.ShapeRange.Count -gt 0) {
        $shp = $hdr.Range.ShapeRange.Item(1)
        if ($shp.Name -eq "image1.jpg") {
            $shp.Name = "image2.jpg"
        }
    }
}

Write-Output "Renamed header/footer logo pictures."
